{"js": "// Apply the UC012 \"Listar Liquida\u00e7\u00f5es Pendentes\" revision:\n// - bump revision history row (version / change-type / date)\n// - small Portuguese copy-edits (accents, punctuation, wording) in the\n//   \"Listar Liquida\u00e7\u00f5es Pendentes\" flow description and its precondition.\n\nasync function replaceOnce(context, searchText, replacement, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = context.document.body.search(searchText, opts);\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. Revision table: version number 1.0.1 -> 1.2.5\nawait replaceOnce(context, \"1.0.1\", \"1.2.5\");\n\n// 2. Revision table: change type Creation -> Update\nawait replaceOnce(context, \"Creation\", \"Update\");\n\n// 3. Revision table: date 04/05/2023 -> 31/05/2023\nawait replaceOnce(context, \"04/05/2023\", \"31/05/2023\");\n\n// 4. Precondition: fix accent on \"usu\u00e1rio\" and add trailing period.\nawait replaceOnce(\n  context,\n  \"O usuario devidamente autenticado e na tela inicial do sistema\",\n  \"O usu\u00e1rio devidamente autenticado e na tela inicial do sistema.\"\n);\n\n// 5. Main flow step 4: accent fixes (\"n\u00famero\", \"di\u00e1ria\", \"di\u00e1rias\").\nawait replaceOnce(context, \"numero da diaria\", \"n\u00famero da di\u00e1ria\");\nawait replaceOnce(context, \"lista de diarias tamb\u00e9m\", \"lista de di\u00e1rias tamb\u00e9m\");\n\n// 6. AF[1] step 2: add trailing period.\nawait replaceOnce(\n  context,\n  \"2. System Apresenta a tela de Detalhar Di\u00e1rias \",\n  \"2. System Apresenta a tela de Detalhar Di\u00e1rias. \"\n);\n\n// 7. AF[2] step 2: remove redundant \"o nome\" after \"onde\".\nawait replaceOnce(\n  context,\n  \"onde o nome dever\u00e1 constar o nome do usu\u00e1rio logado\",\n  \"onde dever\u00e1 constar o nome do usu\u00e1rio logado\"\n);\n\n// 8. AF[3] step 2: add trailing period.\nawait replaceOnce(\n  context,\n  \"2. System Apresenta a tela de Registrar Liquida\u00e7\u00f5es \",\n  \"2. System Apresenta a tela de Registrar Liquida\u00e7\u00f5es. \"\n);\n", "ps1": "# Apply the UC012 \"Listar Liquida\u00e7\u00f5es Pendentes\" revision:\n# - bump revision history row (version / change-type / date)\n# - small Portuguese copy-edits (accents, punctuation, wording) in the\n#   \"Listar Liquida\u00e7\u00f5es Pendentes\" flow description and its precondition.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1. Revision table: version number 1.0.1 -> 1.2.5\nReplace-Text \"1.0.1\" \"1.2.5\"\n\n# 2. Revision table: change type Creation -> Update\nReplace-Text \"Creation\" \"Update\"\n\n# 3. Revision table: date 04/05/2023 -> 31/05/2023\nReplace-Text \"04/05/2023\" \"31/05/2023\"\n\n# 4. Precondition: fix accent on \"usu\u00e1rio\" and add trailing period.\nReplace-Text \"O usuario devidamente autenticado e na tela inicial do sistema\" \"O usu\u00e1rio devidamente autenticado e na tela inicial do sistema.\"\n\n# 5. Main flow step 4: accent fixes (\"n\u00famero\", \"di\u00e1ria\", \"di\u00e1rias\").\nReplace-Text \"numero da diaria\" \"n\u00famero da di\u00e1ria\"\nReplace-Text \"lista de diarias tamb\u00e9m\" \"lista de di\u00e1rias tamb\u00e9m\"\n\n# 6. AF[1] step 2: add trailing period.\nReplace-Text \"2. System Apresenta a tela de Detalhar Di\u00e1rias \" \"2. System Apresenta a tela de Detalhar Di\u00e1rias. \"\n\n# 7. AF[2] step 2: remove redundant \"o nome\" after \"onde\".\nReplace-Text \"onde o nome dever\u00e1 constar o nome do usu\u00e1rio logado\" \"onde dever\u00e1 constar o nome do usu\u00e1rio logado\"\n\n# 8. AF[3] step 2: add trailing period.\nReplace-Text \"2. System Apresenta a tela de Registrar Liquida\u00e7\u00f5es \" \"2. System Apresenta a tela de Registrar Liquida\u00e7\u00f5es. \"\n"}
